# Fruta / hortaliza, semanal
# The weekly data rows (2-19) had their date/volume/price/origin figures
# reshuffled between rows while the descriptive columns (market, region,
# product, category, etc.) stay the same. Re-apply the new per-row values
# for columns D (Fecha), M (Volumen), N (Precio minimo), O (Precio maximo),
# P (Precio promedio ponderado), R (Origen) and S (Precio $/Kg).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($Sheet, $Row, $Fecha, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $Origen, $PrecioKg) {
    $Sheet.Range("D$Row").Value = $Fecha
    $Sheet.Range("M$Row").Value = $Volumen
    $Sheet.Range("N$Row").Value = $PrecioMin
    $Sheet.Range("O$Row").Value = $PrecioMax
    $Sheet.Range("P$Row").Value = $PrecioProm
    $Sheet.Range("R$Row").Value = $Origen
    $Sheet.Range("S$Row").Value = $PrecioKg
}

Set-Row $ws 2  44960 40  7000 7000 7000 "Provincia de Curicó"  3500
Set-Row $ws 3  45001 66  7500 8000 7773 "Provincia de Curicó"  3886
Set-Row $ws 4  44589 60  6000 6000 6000 "Provincia de Curicó"  3000
Set-Row $ws 5  44214 48  6000 6000 6000 "Provincia de Linares" 3000
Set-Row $ws 6  44209 58  6000 6000 6000 "Provincia de Curicó"  3000
Set-Row $ws 7  44592 30  8000 8000 8000 "Provincia de Linares" 4000
Set-Row $ws 8  44211 45  6000 6000 6000 "Provincia de Curicó"  3000
Set-Row $ws 9  44582 150 6000 6500 6233 "Provincia de Curicó"  3116
Set-Row $ws 10 44974 130 7000 7500 7269 "Provincia de Curicó"  3634
Set-Row $ws 11 44586 80  7000 7000 7000 "Provincia de Curicó"  3500
Set-Row $ws 12 44585 160 6500 7000 6750 "Provincia de Curicó"  3375
Set-Row $ws 13 44588 160 6500 7000 6750 "Provincia de Curicó"  3375
Set-Row $ws 14 44628 40  6000 6000 6000 "Provincia de Linares" 3000
Set-Row $ws 15 44606 45  7000 7000 7000 "Provincia de Linares" 3500
Set-Row $ws 16 44587 165 6500 7000 6742 "Provincia de Linares" 3371
Set-Row $ws 17 44627 45  6000 6000 6000 "Provincia de Linares" 3000
Set-Row $ws 18 44959 40  7000 7000 7000 "Provincia de Curicó"  3500
Set-Row $ws 19 44614 45  6000 6000 6000 "Provincia de Linares" 3000
